$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 17:39"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1266785
$ws.Range("C4").Value = 3693
$ws.Range("D4").Value = 213138
$ws.Range("E4").Value = 978685
$ws.Range("F4").Value = 15827
$ws.Range("G4").Value = 163
$ws.Range("H4").Value = 74962

# Row 10: Alemania
$ws.Range("A10").Value = "Alemania"
$ws.Range("B10").Value = 168655
$ws.Range("C10").Value = 493
$ws.Range("D10").Value = 139900
$ws.Range("E10").Value = 21433
$ws.Range("F10").Value = 1884
$ws.Range("G10").Value = 47
$ws.Range("H10").Value = 7322

# Row 25: Pakistan
$ws.Range("A25").Value = "Pakistan"
$ws.Range("B25").Value = 24644
$ws.Range("C25").Value = 1430
$ws.Range("D25").Value = 6464
$ws.Range("E25").Value = 17595
$ws.Range("F25").Value = 111
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 585

# Row 26: Suecia
$ws.Range("A26").Value = "Suecia"
$ws.Range("B26").Value = 24623
$ws.Range("C26").Value = 705
$ws.Range("D26").Value = 4074
$ws.Range("E26").Value = 17509
$ws.Range("F26").Value = 425
$ws.Range("G26").Value = 99
$ws.Range("H26").Value = 3040

# Row 27: Chile
$ws.Range("A27").Value = "Chile"
$ws.Range("B27").Value = 24581
$ws.Range("C27").Value = 1533
$ws.Range("D27").Value = 11664
$ws.Range("E27").Value = 12632
$ws.Range("F27").Value = 493
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 285

# Row 29: Singapur
$ws.Range("A29").Value = "Singapur"
$ws.Range("B29").Value = 20939
$ws.Range("C29").Value = 741
$ws.Range("D29").Value = 1712
$ws.Range("E29").Value = 19207
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 20

# Row 79: Bulgaria
$ws.Range("A79").Value = "Bulgaria"
$ws.Range("B79").Value = 1829
$ws.Range("C79").Value = 51
$ws.Range("D79").Value = 384
$ws.Range("E79").Value = 1361
$ws.Range("F79").Value = 43
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 84

# Row 81: Cuba
$ws.Range("A81").Value = "Cuba"
$ws.Range("B81").Value = 1729
$ws.Range("C81").Value = 26
$ws.Range("D81").Value = 1031
$ws.Range("E81").Value = 625
$ws.Range("F81").Value = 8
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 73

# Row 82: Estonia
$ws.Range("A82").Value = "Estonia"
$ws.Range("B82").Value = 1720
$ws.Range("C82").Value = 7
$ws.Range("D82").Value = 273
$ws.Range("E82").Value = 1391
$ws.Range("F82").Value = 4
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 56

# Row 97: Republica de Chipre
$ws.Range("A97").Value = "Republica de Chipre"
$ws.Range("B97").Value = 889
$ws.Range("C97").Value = 6
$ws.Range("D97").Value = 296
$ws.Range("E97").Value = 578
$ws.Range("F97").Value = 15
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 15

# Row 114: San Marino
$ws.Range("A114").Value = "San Marino"
$ws.Range("B114").Value = 622
$ws.Range("C114").Value = 14
$ws.Range("D114").Value = 106
$ws.Range("E114").Value = 475
$ws.Range("F114").Value = 4
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 41

# Row 115: Georgia
$ws.Range("A115").Value = "Georgia"
$ws.Range("B115").Value = 615
$ws.Range("C115").Value = 5
$ws.Range("D115").Value = 275
$ws.Range("E115").Value = 331
$ws.Range("F115").Value = 6
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 9

# Row 118: Jordania
$ws.Range("A118").Value = "Jordania"
$ws.Range("B118").Value = 484
$ws.Range("C118").Value = 11
$ws.Range("D118").Value = 381
$ws.Range("E118").Value = 94
$ws.Range("F118").Value = 5
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 9

# Row 119: Tanzania
$ws.Range("A119").Value = "Tanzania"
$ws.Range("B119").Value = 480
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 167
$ws.Range("E119").Value = 297
$ws.Range("F119").Value = 7
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 16

# Row 120: Jamaica
$ws.Range("A120").Value = "Jamaica"
$ws.Range("B120").Value = 478
$ws.Range("C120").Value = 5
$ws.Range("D120").Value = 57
$ws.Range("E120").Value = 412
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 9

# Row 121: Guinea-Bisau
$ws.Range("A121").Value = "Guinea-Bisau"
$ws.Range("B121").Value = 475
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 24
$ws.Range("E121").Value = 449
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 2

# Row 186: Botsuana
$ws.Range("A186").Value = "Botsuana"
$ws.Range("B186").Value = 23
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 9
$ws.Range("E186").Value = 13
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 1
